# Update Stundenzettel: add a new "JR" column (with its own "Soll" target
# formula) between the existing "MO" and "KB" columns. The previous "KB"
# column (E) slides one column to the right (F) together with its values
# and formatting; the trailing "Legende" block in columns G/H is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the existing "KB" column (E) one column right, into F -------
# Header (E2 "KB")
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# "Soll" row + monthly data rows (E4:E14)
$ws.Range("E4:E14").Copy()
$ws.Range("F4:F14").PasteSpecial(-4163)
$ws.Range("E4:E14").Copy()
$ws.Range("F4:F14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "noch offen" totals row: move the old formula/format from E19 to F19 and
# repoint it at the (now shifted) column F data.
$ws.Range("E19").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F19").Formula = "=F4-SUM(F5:F14)"

# Old placeholder cells below the data table (F15/F16) are no longer needed
# once the real content has moved in above them.
$ws.Range("F15").Clear()
$ws.Range("F16").Clear()

# --- 2. Populate the new "JR" column (E) -----------------------------------
$ws.Range("E2").Value = "JR"
$ws.Range("E4").Formula = "=(43/2)+3*43"

$ws.Range("E5:E12").ClearFormats()
$ws.Range("E5:E12").Value = "-"

$ws.Range("E19").ClearFormats()
$ws.Range("E19").Formula = "=E4-E15-E16-E17"

# --- 3. New "Soll bis 31.3.23" label under the header row ------------------
$ws.Range("A3").Value = "Soll bis 31.3.23"

# --- 4. Restore the active selection ---------------------------------------
$ws.Range("A4").Select()
